$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# Columns D (4) and E (5): stored width 15.7109375 -> 14.7109375
# Column L (12): stored width 16 -> 10
# Note: the ColumnWidth COM property is expressed in "characters" and gets
# rounded to the nearest pixel internally, so we pick the ColumnWidth value
# whose resulting stored width is closest to the target.
$ws.Columns.Item(4).ColumnWidth = 13.75
$ws.Columns.Item(5).ColumnWidth = 13.75
$ws.Columns.Item(12).ColumnWidth = 9.166666666666666

# --- Text label updates (row 3: L3/M3) ---
$ws.Range("L3").Value = "small"
$ws.Range("M3").Value = "**"

# --- Updated statistics values (rows 2-4) ---
# Row 2
$ws.Range("D2").Value = [double]"0.011728814653501311"
$ws.Range("E2").Value = [double]"0.011728814653501285"
$ws.Range("F2").Value = [double]"-0.0036430571514536325"
$ws.Range("G2").Value = [double]"-0.051399759494407185"
$ws.Range("H2").Value = [double]"6.3993123600096151"
$ws.Range("J2").Value = [double]"490"
$ws.Range("K2").Value = [double]"0.012891460968359612"

# Row 3
$ws.Range("D3").Value = [double]"0.0012944081158074961"
$ws.Range("E3").Value = [double]"0.0012944081158074994"
$ws.Range("F3").Value = [double]"0.0023769071387462759"
$ws.Range("G3").Value = [double]"0.037166259800806629"
$ws.Range("H3").Value = [double]"10.471349452636249"
$ws.Range("J3").Value = [double]"490"
$ws.Range("K3").Value = [double]"0.020922974839795976"

# Row 4
$ws.Range("D4").Value = [double]"7.4549944846953905e-10"
$ws.Range("E4").Value = [double]"1.4909988799161056e-09"
$ws.Range("F4").Value = [double]"-0.0096630214416535409"
$ws.Range("G4").Value = [double]"-0.12420256015361938"
$ws.Range("H4").Value = [double]"39.443393819986085"
$ws.Range("J4").Value = [double]"490"
$ws.Range("K4").Value = [double]"0.074499737423104145"
